$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.176.09"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "3.404.91"
$ws.Range("E3").Value = "  +3.02%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.00%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.63"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +3.47%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.30"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +5.41%  "

$ws.Range("E7").Value = "  +4.56%  "

$ws.Range("D8").Value = "3.395.22"
$ws.Range("E8").Value = "  +3.32%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +14.20%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.637"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +5.20%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.21"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +2.76%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +7.05%  "

$ws.Range("E14").Value = "  +4.56%  "

$ws.Range("D15").Value = "3.936.58"
$ws.Range("E15").Value = "  +2.41%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.39"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +5.57%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.402.86"
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +2.06%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +4.66%  "

$ws.Range("D20").Value = "65.176.35"
$ws.Range("E20").Value = "  +3.29%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +4.23%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.98"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +17.72%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.01"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +16.05%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.16"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +4.17%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.52"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +6.58%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.59"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +3.97%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +2.66%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.88"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +6.68%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.93"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +4.86%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.69"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +6.84%  "

$ws.Range("E31").Value = "  +3.76%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.59"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +3.93%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "583.66"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("E34").Value = "  +5.16%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.29"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +4.50%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -3.81%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.22"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").Value = "0.0₃0767"
$ws.Range("E39").Value = "  +6.03%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +1.00%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.376"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +4.07%  "

$ws.Range("D42").Value = "3.122.71"
$ws.Range("E42").Value = "  +0.51%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +0.02%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.89"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +4.91%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +4.07%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +4.84%  "

$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("E48").Value = "  +5.29%  "

$ws.Range("E49").Value = "  -0.17%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.50"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +7.17%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.70"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +4.82%  "
